$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 46039
$ws.Range("B2").Value = 135.59
$ws.Range("C2").Value = 126.8
$ws.Range("D2").Value = 122.91
$ws.Range("E2").Value = 119.84
$ws.Range("F2").Value = 119.14
$ws.Range("G2").Value = 118.84
$ws.Range("H2").Value = 120.35
$ws.Range("I2").Value = 123.12
$ws.Range("J2").Value = 126.82
$ws.Range("K2").Value = 121.98
$ws.Range("L2").Value = 118.26
$ws.Range("M2").Value = 117.99
$ws.Range("N2").Value = 116.7
$ws.Range("O2").Value = 112.07
$ws.Range("P2").Value = 115.82
$ws.Range("Q2").Value = 120.27
$ws.Range("R2").Value = 124.36
$ws.Range("S2").Value = 135.98
$ws.Range("T2").Value = 142.85
$ws.Range("U2").Value = 140.91
$ws.Range("V2").Value = 155.23
$ws.Range("W2").Value = 154.7
$ws.Range("X2").Value = 130.11
$ws.Range("Y2").Value = 133.09
$ws.Range("Z2").Value = 127.24
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 143.28
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 154.96
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 141.88
$ws.Range("AG2").Value = "1h-16h"
